# Give extra decimal precision to five percentage figures in the
# "Source Type: SME Associations (Most Widely Used)" table on the
# Summary sheet:
#   C32: 43.7  -> 43.74
#   D32: 62.2  -> 62.25
#   B34: 21.6  -> 21.56
#   D34: 26.3  -> 26.25
#   B36: 81.7  -> 81.71
#
# These figures are stored as text (not numbers) in the workbook. Excel
# auto-converts a numeric-looking string typed into a cell into a real
# number, so for each cell we briefly force a text number-format before
# writing the new value, then restore the cell's original style so no
# other formatting is changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$rng = $ws.Range("C32")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "43.74"
$rng.Style = $origStyle

$rng = $ws.Range("D32")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "62.25"
$rng.Style = $origStyle

$rng = $ws.Range("B34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "21.56"
$rng.Style = $origStyle

$rng = $ws.Range("D34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "26.25"
$rng.Style = $origStyle

$rng = $ws.Range("B36")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "81.71"
$rng.Style = $origStyle
